# Update the "Website / Source" column (D) so that each cell shows the
# actual hyperlink URL text instead of repeating the tool's name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = "https://datakitchen.io/"
$ws.Range("D4").Value  = "https://airbyte.com/"
$ws.Range("D5").Value  = "https://www.fivetran.com/"
$ws.Range("D6").Value  = "https://dagster.io/"
$ws.Range("D7").Value  = "https://www.prefect.io/"
$ws.Range("D8").Value  = "https://www.matillion.com/"
$ws.Range("D9").Value  = "https://www.talend.com/"
$ws.Range("D10").Value = "https://www.getdbt.com/"
$ws.Range("D11").Value = "https://www.ibm.com/products/streamsets"
$ws.Range("D12").Value = "https://nifi.apache.org/"
$ws.Range("D13").Value = "https://hevodata.com/"
$ws.Range("D14").Value = "https://www.snaplogic.com/"
$ws.Range("D15").Value = "https://www.informatica.com/"
$ws.Range("D16").Value = "https://admin.google.com/ServiceNotAllowed?application=995920231026&source=scrip&continue=https://cloud.google.com/dataflow"
$ws.Range("D17").Value = "https://aws.amazon.com/glue/"
$ws.Range("D18").Value = "https://azure.microsoft.com/en-us/products/data-factory/"
$ws.Range("D19").Value = "https://sourceforge.net/projects/pentaho/"
$ws.Range("D20").Value = "https://www.alteryx.com/about-us/trifacta-is-now-alteryx-designer-cloud"
$ws.Range("D21").Value = "https://airflow.apache.org/"
$ws.Range("D22").Value = "https://greatexpectations.io/"

# Restore the view to show the top of the sheet with F22 selected instead
# of the previous scroll position / selection (C24).
$null = $ws.Range("F22").Select()
